$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.232.06"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.45%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.432.47"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.18%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.58"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.43"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.46%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.514"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.95%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.499"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.35%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.19"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.13%  "

$ws.Range("E11").Value = "  +1.34%  "

$ws.Range("E12").Value = "  +2.49%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.72"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.38%  "

$ws.Range("E14").Value = "  +2.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.806.60"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.02%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.411.25"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.47%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.833"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +3.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.173.74"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.36%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.33"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.81%  "

$ws.Range("E20").Value = "  +1.33%  "

$ws.Range("E21").Value = "  +1.82%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.63"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.39"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.06%  "

$ws.Range("E24").Value = "  +3.57%  "

$ws.Range("E25").Value = "  +1.34%  "

$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("E27").Value = "  +2.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.21"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -6.69%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.51"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.91"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.12%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.120"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +16.98%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.64"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +8.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.17"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.51%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0767"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.93"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.64%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.56"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +5.67%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "130.84"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +20.70%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.94"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +5.09%  "

$ws.Range("E40").Value = "  -0.76%  "

$ws.Range("E41").Value = "  +0.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.01"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -6.11%  "

$ws.Range("E43").Value = "  +2.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.957.71"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.18"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.13%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.88"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.54%  "

$ws.Range("E47").Value = "  +0.76%  "

$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.65"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +8.92%  "

$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.647.54"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.60%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.52"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.33%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.57"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.06%  "
